$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target word list for rows 63..1028 (shared-string index realignment
# after removing the "cada" and "través" entries from sharedStrings.xml)
$words = @(
  'cultura',
  'apoyar',
  'difíciles',
  'laboral',
  'curso',
  'ayudando',
  'radio',
  'ahora',
  'empezar',
  'escritura',
  'puedes',
  'mundo',
  'gracias',
  'fondos',
  'adultos',
  'espacio',
  'escolar',
  'latinas',
  'programa',
  'mano',
  'hacemos',
  'sueño',
  'hombres',
  'puede',
  'áreas',
  'seguir',
  'mejor',
  'corazón',
  'oportunidad',
  'alcoholismo',
  'menos',
  'habla',
  'becas',
  'comunidad.',
  'herramientas',
  'lograr',
  'dando',
  'salvador',
  'deseamos',
  'valores',
  'agentes',
  'pasa',
  'asistencia',
  'laborales',
  'mamas',
  'práctica',
  'instituto',
  '30',
  'dedica',
  'hospital',
  'objetivo',
  'trabajo',
  'dios',
  'autónoma',
  'cuerpos',
  'maestro',
  'hoy',
  'indígena',
  'natural',
  'donación',
  'expresión',
  'romper',
  'galveston',
  'recaudando',
  'recuperación',
  'forma',
  'sido',
  'especial',
  'apoya',
  'estaciones',
  'nunca',
  'red',
  'preparación',
  'amigos',
  'realizar',
  'elecciones',
  'esperanza',
  'dentales',
  'realidad',
  'conciencia',
  'habilidades',
  'alegría',
  'trabajadores',
  'latinx',
  'latina',
  'latino',
  'manos',
  'dos',
  'bajos',
  'equipos',
  'local',
  'papel',
  'guatemala',
  'necesitan',
  'interés',
  'cambiar',
  'buscamos',
  'luz',
  'mantenerlas',
  'informadas.',
  'públicas',
  'alcanzar',
  'méxico.',
  'aquellos',
  'opciones',
  'fundación',
  'actualmente',
  'casa',
  'generar',
  'juegan',
  'voz',
  'basurero',
  'económico',
  'fortaleciendo',
  'conectar',
  'tiempo',
  'año',
  'apoyemos',
  'glorietas',
  'respeto',
  'parar',
  '4',
  'ambulancia',
  'mujer',
  'viven',
  'desarrollo',
  'niño',
  'nivel',
  'sustentar',
  'transformamos',
  'pasos',
  'mayores',
  'nueva',
  'básicas',
  'prevenir',
  'trata',
  'liderazgo',
  'empresa',
  'mas',
  'quién',
  'organización',
  'proteger',
  'juntos',
  'derecho',
  'sexual',
  'cambio',
  'shriners',
  'vulnerables',
  'jovenes',
  'posible',
  'rural',
  'creativa',
  'aprovechemos',
  'circundantes',
  'ciudad',
  'futuro.',
  'inquietud',
  'promueve',
  '15,000',
  'pobreza',
  'presenta',
  'sucede.',
  'tx',
  'viajar',
  'vez',
  'igual',
  'adicciones',
  'humildes',
  'mortalidad',
  'vel',
  'educar',
  'superar',
  'queridos',
  'sororidad',
  'hispanos',
  'aporta',
  'esfuerzo',
  'afecten',
  'persona',
  'educativa',
  'evento',
  'alumnas',
  'inclusión',
  'supervisada',
  'casas',
  'comunitaria',
  'propio',
  'capacitaciones',
  'latinavida',
  'padres.',
  'adquisición',
  'stamos',
  'cultivo',
  'consentimiento',
  'incubadora',
  'pobreza.',
  'sirviendo',
  'diseñar',
  'navidad',
  'alfa',
  'trajes',
  'vandalismo',
  'paso',
  'capacitación',
  'cuentan',
  'elegimos',
  'defender',
  'ciberespacio',
  'participar',
  'sue',
  'reciban',
  'teóricos',
  'trabajo.',
  '8m',
  'acceso',
  'ntilde',
  'escuelas',
  'proporcionados',
  'huracán',
  'ejercicios',
  'donaciones',
  'drogadicción',
  'honor.',
  'sitúan',
  'denunció',
  'estatus',
  'avanzando',
  'igualdad',
  'colaborar',
  'dia.',
  'enfoca',
  'activo',
  'convenio',
  'capacitarte',
  'reflexión',
  'empoderando',
  'apoye',
  '100',
  'extrema',
  'aprendizaje',
  'anunciemos',
  'exactamente',
  'dólar',
  'adolescente',
  'transplante.',
  'alta',
  'argentina.',
  'acción',
  'dls.',
  'unesco.',
  'gloria',
  'financiera',
  'niñez',
  'víveres',
  'fortaleza…',
  'app',
  'infantil',
  'escuchar',
  'creación',
  'probablemente',
  'pobre.',
  'atender',
  'seriamos',
  'sanación',
  'administración',
  'fracaso',
  'mejora',
  'promoviendo',
  'nuevas',
  'participantes',
  'refugio',
  'explotación',
  'negativa',
  'etc',
  'movilizamos',
  'trabaja',
  'todas',
  'campesinos',
  'encontramos',
  '6',
  'ofrecerles',
  'social.',
  'guerreros',
  'críticas',
  'tx.',
  'continuar',
  'sabes',
  'duper',
  'etnias',
  'mejor.',
  'impulsa',
  'comunitaria.',
  'condiciones',
  '97',
  'integrales',
  'valoración',
  '¡aprende',
  'biocombustibles',
  'universitarios',
  'tatuajes',
  'proceso',
  '¡manténlo',
  'pasar.',
  'heridas',
  'salvarle',
  'trump',
  'estatal.',
  'construyen',
  'justo',
  'afectiva',
  'voces',
  'sentirse',
  'perder',
  'derivados',
  'cocinar',
  'pone',
  'obtendrán',
  'bolsas',
  'impacto',
  'viejos',
  'área',
  'seguridad',
  'incendios.',
  '1',
  'autosuficiente',
  'maravillas',
  'entregando',
  'llega',
  'vejez',
  'simulación',
  'estratégico',
  'vestido',
  'indígenas.',
  'reconstruir',
  'falta',
  'especializada.',
  'cambiaremos',
  'comparte',
  'inspirando',
  'cuesta',
  'especializada',
  'azotaron',
  'damos',
  'inspirar',
  'utilizaremos',
  'aviones',
  'q',
  'si',
  'desventajas',
  'después',
  'cuidados',
  'primero',
  'salvarlos',
  'manera',
  'ayudanos',
  'edad.',
  'atención',
  'comunitario',
  'víctimas',
  'procesos',
  'diferentes',
  'penales',
  'cólera',
  'mayor',
  'literario',
  'baile',
  'indocumentadas',
  'personal',
  'gente',
  'nace',
  'surgeons',
  'condado',
  'estudios',
  'necesidades',
  'requiere',
  'necesitadas',
  'fortaleza',
  '2019',
  'generacion',
  'tiempos',
  'puerto',
  'fines',
  'negocio…',
  'proporcionar',
  '¡hagamos',
  'guajira',
  'protección',
  'practica',
  'abandone',
  'miedo',
  'afrontar',
  'vuelo',
  'dar',
  'educando',
  'trabajamos',
  'indocumentad',
  'instituciones',
  'travez',
  'sepan',
  'mérida',
  'ayudarnos',
  '25',
  'multas',
  'muchas',
  'emocionales',
  'histórico',
  'abogacía',
  'cubrir',
  'baltimore.',
  'matthew.',
  'encestale',
  'pertenecientes',
  'glasswing',
  'creativa.',
  'dia',
  'recibir',
  'toda',
  'creemos',
  'vivamus',
  'migración',
  'superior.',
  'dificultades',
  'odontológica',
  '¡amplifique',
  'gol',
  'hogar',
  'arreglarlo',
  'preventivas',
  'presidenciales',
  'remodelación',
  'empoderamiento',
  'producción',
  'asista',
  '7',
  'retórica',
  'lectura',
  'textil',
  'undólarcadadía',
  'latinavida360',
  'aplicar',
  'algas',
  'pequeña',
  'competencias',
  'silencio',
  'asalto',
  'matlapa',
  'profesional',
  'obra',
  'fuente',
  'autóctonas',
  'piloto',
  'apoderan',
  'cuidado',
  'congénita',
  'brindara',
  'líderes',
  'universidad',
  'resultados',
  'podemos',
  'erradicación',
  'extinción',
  'prestando',
  'años.',
  'resguardo',
  'salir',
  'detenciones',
  'echo',
  'dona',
  'potencial',
  'colombia',
  'reforestar',
  'etnia',
  'golpeado',
  'primeros',
  'internet.',
  'principios',
  'corporal',
  'gratuitamente',
  'sobrevivientes.',
  'enciende',
  'inicial',
  'lucharemos',
  'fundacion',
  'legislativo',
  'riesgo',
  'oncología',
  'financieros',
  'luchar',
  'guadalupe',
  'miembros',
  'permitamos',
  'cambio.',
  'guatemaltecas',
  'estudiantes.',
  'ayúdame',
  '¡lo',
  'temprana',
  'conocimientos',
  'doctora',
  'segundo',
  'llevárlos',
  'segunda',
  'importate',
  'contribuir',
  'brindar',
  'momento',
  '¡la',
  'hispana',
  'requerimos',
  'digo',
  'bajo',
  'relaciones',
  'tan',
  'asegurarnos',
  '11',
  '10',
  'sensibilización',
  'número',
  '2015',
  'evitar',
  'ayúdenos',
  'historia.',
  'motriz',
  'project',
  'escenario',
  'llevar',
  'sed',
  'territorios.',
  'colombiano',
  'juguetes',
  'san',
  'cooperación',
  'bailar',
  'salvar',
  'costo',
  'emociones..',
  'vulnerable',
  'importe',
  'dentro',
  'profesionales',
  'blanca',
  'buentrato',
  'peligro',
  'sólo',
  'parece',
  'comunitarios',
  'méxico',
  'derek',
  '¡tenemos',
  'vulnerabilidad',
  'avecinan',
  'clase',
  'lucro',
  'centro',
  '1ro',
  'limitados',
  'estimula',
  'visa',
  'éxito',
  'amar',
  'pais',
  'unidad',
  'animales',
  'renal',
  'digna',
  'médicos',
  'be-you-tiful',
  'reflejan',
  'organizaciones',
  'metas',
  'muxer',
  '30,000',
  'empoderados',
  'camioneta',
  'laoreet',
  'represión',
  'mantenimiento',
  'elaborar',
  'conformación',
  'conocimiento',
  'saludable',
  'educamos',
  'sexualidad',
  'hace',
  'localidades',
  'realizado',
  'hacerlo',
  'adolescencia',
  'reduzcarmos',
  'emprende',
  '3',
  'internacional',
  'rico',
  'declarada',
  'quemaduras',
  'formas',
  'hope',
  'dale',
  'toco',
  'dignidad',
  'últimas',
  'familiar',
  '12',
  'estafadores',
  'públicos',
  'alfabetización',
  'aportando',
  'tratamiento',
  'votantes',
  'individuo',
  'literaria',
  'registramos',
  'lugar',
  'curativas',
  'unidos',
  'aproximadamente',
  'cirugías',
  'mision',
  'marinas',
  'ganas',
  'surgen',
  '//unidadlatinanj.org/get-involved/donate/',
  'segundos',
  'médica',
  'alimentación',
  'traerá',
  'kantaya',
  'wayuú',
  'momentos',
  'perdido',
  'guerra',
  'méxicanos',
  'entrena',
  'humano',
  '038',
  'individual',
  'discriminadas.',
  'l.o.v.e',
  'conferencia',
  'hacía',
  'purpura',
  '¡ayúdenos',
  'progreso',
  'contenido',
  'voluntarios',
  'juego',
  'mexicanas.',
  'emprendimiento',
  'nutrición',
  'import',
  'alguien',
  'lideres',
  'local.',
  'amplio',
  'frente',
  'asociacion',
  'explotaba',
  'avanzar',
  'experiencia',
  'futuro',
  'espació',
  'ambiente',
  'haití',
  'experimentas',
  'legislatura',
  'ofrece',
  'gratis.',
  'reduciendo',
  'tocan',
  'concluir',
  'dolor',
  'dominicanas.',
  'desastre',
  'agente',
  'previniendo',
  'muebles',
  'simple',
  'empresarias',
  'órganos',
  'misericordia',
  'lucha',
  'tormenta',
  'fronteras.',
  'medios',
  'señor',
  'desarrollarse',
  'económica',
  'nacional',
  'interna',
  'perdura.',
  'doméstica',
  'dui',
  'cualquiera',
  'económicas.',
  'transformarnos',
  'fiestas',
  'amiga',
  'enriquecimiento',
  'pa',
  'arte',
  'ahorros',
  'obstáculo',
  'orioles',
  'cena',
  'cuidar',
  'universal',
  'dame',
  'enseñándoles',
  'alto',
  'reducir',
  'día',
  'ayudemos',
  'usado',
  'intelectual',
  'sanen',
  'pie',
  'donec',
  'extracurricular',
  'quemados',
  'pasaporte',
  'educacion',
  'desconstruir',
  'proyectos',
  'reflexiona',
  'urgencia',
  'planeta.',
  'orientada',
  'evitamos',
  'situaciones',
  'drogadicion',
  'sostenibilidad',
  'buen',
  'premiar',
  'adentrar',
  'servicios',
  'zonas',
  'encuentran',
  'cáncer',
  'comenzado',
  'anil',
  'iv',
  'peruanos',
  'sagittis',
  'busca',
  'comidagratis',
  'relacionamiento',
  'alimentar',
  'establecer',
  'medio',
  'media',
  'juan',
  'asegurarse',
  'quehacer',
  'necesarios',
  'ninas',
  'combatir',
  'provienen',
  'http',
  'junto',
  'debería',
  'solidarízate',
  'cuya',
  'poblana',
  'desiguales',
  'padres',
  'merecemos',
  'deportaciones',
  'tierra',
  'haga',
  'dado',
  'aún',
  'estatal',
  'tasa',
  'injusticias',
  'navideña',
  'sumas',
  'talleres',
  'representantes',
  'proveerlos',
  'lacus',
  'sabiduría.',
  'medias-becas',
  'limite',
  'remoción',
  'desarrollar',
  'respiro',
  'superacion',
  'sonrisas',
  'sonoma',
  'discriminación',
  'protegiendo',
  'terremotos',
  'racismo',
  'salvador.',
  'laboratorio',
  'invitamos',
  'temporales',
  'avión',
  'convocado',
  'atrapados',
  'transformacion',
  'resistiendo',
  'maria',
  'fortalecer',
  'maquillaje',
  'desean',
  'marca',
  'enajo',
  'transportar',
  'autonomía',
  '260',
  'unirse',
  'mixteca',
  'faucibus',
  'menstruación',
  'sostenibles',
  'américa',
  'augue',
  'realidades',
  'máximo',
  'abrigo',
  'trato',
  'unirnos',
  'redes',
  'acompáñanos',
  'ustedes',
  'yachay',
  'semanal',
  'permitir',
  'luchando',
  'recursos.',
  'york.',
  'difusión',
  'innovadores',
  'accesorios',
  'contactos',
  'patrimonio',
  'inquilino',
  'equipar',
  'autoempleo',
  'legal',
  '¡ayúdanos',
  'decisiones',
  'formación',
  'reforzar',
  'xenofobia',
  'migrantes',
  'oficio',
  'sinfín',
  'embarazo',
  'pueda',
  'navideñas',
  'stem',
  'nuevo',
  'transmitir',
  'colectiva',
  'emocional',
  'recaudación',
  'perdieron',
  'super',
  'protagonismo',
  'conservacion',
  'materna',
  'construcción',
  'ensenarle',
  'union',
  'envía',
  'secuelas',
  'origina',
  'abogando',
  'venezolanas',
  'promotoras',
  'diversidad',
  'tapachula.',
  'propagación',
  'conocer',
  'encuentro',
  'obtener',
  'ofreciéndole',
  'bueno',
  'alas',
  'rápido',
  'transformar',
  'ayude',
  'desigualdades',
  'ofrecemos',
  'cancer.',
  'esperando',
  'lúdico',
  'mar',
  'comida',
  'homicidios',
  'desvía',
  'humanos',
  'adelante',
  'brindamos',
  'fase',
  'peru',
  'haciendo',
  '700',
  'talento',
  'auctor.',
  'digno.',
  'tramitar',
  'indígenas',
  'circulo',
  'rutrum',
  'bienes',
  'vive',
  'impactados',
  'contribuya',
  'recibido',
  '¡imagínate',
  'únete',
  'marely',
  'cerradas',
  'vivencia',
  'ap',
  'despensa',
  'vidas.',
  'planeta',
  'capacidad',
  'alcance',
  'recupera',
  'empoderadas',
  'recupere',
  'estigmatizadas',
  'disfrutar',
  'ninos',
  'adolecentes',
  'organizando',
  'leer',
  'saber',
  'involucra',
  '¡contribuye',
  '40',
  'asegurar',
  '5',
  'test',
  'intrafamiliar',
  'indígen',
  'pedir',
  'temor',
  'visual.',
  'inclusivo',
  'fundó',
  'verdes',
  'reduceelriesgo',
  'pobres',
  'tomar',
  'transportamos',
  'miles',
  'cocina',
  'convivencia',
  'línea',
  'executivas',
  'ningún',
  'técnicas',
  'empleo.',
  'ortega',
  'crónica',
  'misión',
  'humanidad',
  'oral',
  'encender',
  'colegios'
)

# Row offsets (0-based, relative to row 63) whose target text looks like a
# pure number and must be forced to Text so Excel does not auto-convert it.
$textForce = @(47, 121, 158, 237, 279, 293, 323, 377, 401, 433, 514, 515, 518, 571, 590, 602, 634, 810, 899, 935, 937)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = 63 + $i
    $cell = $ws.Cells.Item($row, 2)
    if ($textForce -contains $i) {
        $cell.NumberFormat = "@"
        $cell.Value = $words[$i]
        $cell.ClearFormats() | Out-Null
    } else {
        $cell.Value = $words[$i]
    }
}

# Word-count decrements
$ws.Cells.Item(92, 3).Value = 2
$ws.Cells.Item(93, 3).Value = 2
$ws.Cells.Item(214, 3).Value = 1
$ws.Cells.Item(215, 3).Value = 1

# Remove the last two data rows (1029 and 1030)
$ws.Rows.Item(1029).Delete() | Out-Null
$ws.Rows.Item(1029).Delete() | Out-Null

Write-Output "edit complete"
